# Generate Report for Handback
# Refresh the handback-status report timestamps / status for the
# "3a096a50-3b66-4e51-8922-996029d52578" entry (row 2/3 of each sheet,
# since rows 2 and 3 previously shared identical values).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-25 06:15:11"
$wsOverview.Range("G3").Value = "2016-08-25 06:15:11"

# --- zh-cn sheet: Priority (E), Correspond Handoff Datetime (H),
#     Correspond Handback DateTime (K) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-25 06:14:59"
$wsZhCn.Range("H3").Value = "2016-08-25 06:14:59"
$wsZhCn.Range("K2").Value = "2016-08-25 06:15:32"
$wsZhCn.Range("K3").Value = "2016-08-25 06:15:32"

# --- de-de sheet: Priority (E), Correspond Handback DateTime (K) ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("K2").Value = "2016-08-25 06:15:40"
$wsDeDe.Range("K3").Value = "2016-08-25 06:15:40"
